$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (shared strings with rich-text runs) ---
# A8: "Volume 32   Number  7" -> "Volume 32   Number  8"
$ws.Range("A8").Characters(21, 1).Text = "8"

# C9: "Report Covering the Week  2/10/2025  Through  2/16/2025"
#  -> "Report Covering the Week  2/17/2025  Through  2/23/2025"
$ws.Range("C9").Characters(27, 9).Text = "2/17/2025"
$ws.Range("C9").Characters(47, 9).Text = "2/23/2025"

# --- Crime-complaint grid (rows 14-31) ---
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("N15").Value2 = -85.714285714285
$ws.Range("F16").Value2 = 2
$ws.Range("G16").Value2 = 7
$ws.Range("H16").Value2 = -71.428571428571
$ws.Range("J16").Value2 = 16
$ws.Range("K16").Value2 = -68.75
$ws.Range("L16").Value2 = -73.684210526315
$ws.Range("N16").Value2 = -96.350364963503
$ws.Range("C17").Value2 = 3
$ws.Range("D17").Value2 = 9
$ws.Range("E17").Value2 = -66.666666666666
$ws.Range("F17").Value2 = 10
$ws.Range("G17").Value2 = 16
$ws.Range("H17").Value2 = -37.5
$ws.Range("I17").Value2 = 21
$ws.Range("J17").Value2 = 23
$ws.Range("K17").Value2 = -8.695652173913
$ws.Range("L17").Value2 = -16
$ws.Range("M17").Value2 = 425
$ws.Range("N17").Value2 = -16
$ws.Range("C18").Value2 = 5
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 66.666666666666
$ws.Range("F18").Value2 = 12
$ws.Range("G18").Value2 = 9
$ws.Range("H18").Value2 = 33.333333333333
$ws.Range("I18").Value2 = 22
$ws.Range("J18").Value2 = 22
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = -50
$ws.Range("M18").Value2 = 22.222222222222
$ws.Range("N18").Value2 = -81.818181818181
$ws.Range("C19").Value2 = 8
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = -11.111111111111
$ws.Range("F19").Value2 = 18
$ws.Range("G19").Value2 = 34
$ws.Range("H19").Value2 = -47.058823529411
$ws.Range("I19").Value2 = 34
$ws.Range("J19").Value2 = 60
$ws.Range("K19").Value2 = -43.333333333333
$ws.Range("L19").Value2 = -65.306122448979
$ws.Range("M19").Value2 = -39.285714285714
$ws.Range("N19").Value2 = -42.372881355932
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = -50
$ws.Range("F20").Value2 = 6
$ws.Range("G20").Value2 = 14
$ws.Range("H20").Value2 = -57.142857142857
$ws.Range("I20").Value2 = 7
$ws.Range("J20").Value2 = 23
$ws.Range("K20").Value2 = -69.565217391304
$ws.Range("L20").Value2 = -58.823529411764
$ws.Range("M20").Value2 = 0
$ws.Range("N20").Value2 = -96.276595744680
$ws.Range("C21").Value2 = 17
$ws.Range("D21").Value2 = 24
$ws.Range("E21").Value2 = -29.166666666666
$ws.Range("F21").Value2 = 48
$ws.Range("G21").Value2 = 80
$ws.Range("H21").Value2 = -40
$ws.Range("I21").Value2 = 90
$ws.Range("J21").Value2 = 146
$ws.Range("K21").Value2 = -38.356164383561
$ws.Range("L21").Value2 = -56.097560975609
$ws.Range("M21").Value2 = -14.285714285714
$ws.Range("N21").Value2 = -83.271375464684
$ws.Range("G22").Value2 = 1
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("F23").Copy($ws.Range("D23"))
$ws.Range("D23").Value2 = 1
$ws.Range("H23").Copy($ws.Range("E23"))
$ws.Range("E23").Value2 = -100
$ws.Range("J23").Value2 = 5
$ws.Range("K23").Value2 = -40
$ws.Range("C24").Value2 = 33
$ws.Range("D24").Value2 = 24
$ws.Range("E24").Value2 = 37.5
$ws.Range("F24").Value2 = 97
$ws.Range("G24").Value2 = 70
$ws.Range("H24").Value2 = 38.571428571428
$ws.Range("I24").Value2 = 195
$ws.Range("J24").Value2 = 145
$ws.Range("K24").Value2 = 34.482758620689
$ws.Range("L24").Value2 = -26.415094339622
$ws.Range("M24").Value2 = 80.555555555555
$ws.Range("C25").Value2 = 28
$ws.Range("D25").Value2 = 17
$ws.Range("E25").Value2 = 64.705882352941
$ws.Range("F25").Value2 = 78
$ws.Range("G25").Value2 = 42
$ws.Range("H25").Value2 = 85.714285714285
$ws.Range("I25").Value2 = 139
$ws.Range("J25").Value2 = 91
$ws.Range("K25").Value2 = 52.747252747252
$ws.Range("L25").Value2 = -36.238532110091
$ws.Range("C26").Value2 = 2
$ws.Range("D26").Value2 = 3
$ws.Range("E26").Value2 = -33.333333333333
$ws.Range("F26").Value2 = 8
$ws.Range("G26").Value2 = 11
$ws.Range("H26").Value2 = -27.272727272727
$ws.Range("I26").Value2 = 31
$ws.Range("J26").Value2 = 24
$ws.Range("K26").Value2 = 29.166666666666
$ws.Range("L26").Value2 = -32.608695652173
$ws.Range("M26").Value2 = 3.333333333333
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("E14").Copy($ws.Range("H27"))
$ws.Range("L27").Value2 = -33.333333333333
$ws.Range("F28").Value2 = 5
$ws.Range("G28").Value2 = 4
$ws.Range("H28").Value2 = 25
$ws.Range("J28").Value2 = 15
$ws.Range("K28").Value2 = -40
$ws.Range("L28").Value2 = 12.5
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))
$ws.Range("C14").Copy($ws.Range("F31"))
